$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column E (copy format from the existing header cells)
$ws.Range("E1").Value = "QUANTIDADE"
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for column E, rows 2-47
$values = @(
    2347,
    227,
    473,
    1186,
    3237,
    7508,
    14967,
    25670,
    39616,
    55609,
    73234,
    90291,
    106553,
    121897,
    134801,
    144874,
    151574,
    155537,
    156735,
    155235,
    150933,
    144627,
    136501,
    125251,
    115133,
    103238,
    91137,
    78786,
    66557,
    54955,
    44420,
    34696,
    25554,
    18461,
    12560,
    8003,
    4785,
    2670,
    1286,
    533,
    227,
    59,
    10,
    4,
    0,
    0
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $values[$i]
}
